$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7785152196884155
$ws.Range("B1").Value = 2.195825576782227
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.904543876647949
$ws.Range("E1").Value = 1.086121439933777
